# chore: update Sheets via scheduled runner
# Refresh cached market-price / profit figures across the leve-profit sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Mirrors a scheduled data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 36328.5
$ws.Range("J3").Value = 36328.5
$ws.Range("L3").Value = 36328.5
$ws.Range("N3").Value = -36556.5
$ws.Range("H19").Value = 351.77777
$ws.Range("I19").Value = 600
$ws.Range("J19").Value = 280.85715
$ws.Range("K19").Value = 600
$ws.Range("L19").Value = 280.85715
$ws.Range("M19").Value = -425
$ws.Range("N19").Value = -630.85715
$ws.Range("H102").Value = 36328.5
$ws.Range("J102").Value = 36328.5
$ws.Range("L102").Value = 36328.5
$ws.Range("N102").Value = -42818.5
$ws.Range("H132").Value = 2674.3035
$ws.Range("I132").Value = 1505.4131
$ws.Range("J132").Value = 8051.2
$ws.Range("K132").Value = 4516.2393
$ws.Range("L132").Value = 24153.6
$ws.Range("M132").Value = -1986.2393
$ws.Range("N132").Value = -29213.6
$ws.Range("H135").Value = 24334
$ws.Range("I135").Value = 26599.385
$ws.Range("J135").Value = 2246.5
$ws.Range("K135").Value = 239394.465
$ws.Range("L135").Value = 20218.5
$ws.Range("M135").Value = -236859.465
$ws.Range("N135").Value = -25288.5
$ws.Range("H137").Value = 3334506.5
$ws.Range("I137").Value = 1563668.9
$ws.Range("J137").Value = 7693491
$ws.Range("K137").Value = 4691006.699999999
$ws.Range("L137").Value = 23080473
$ws.Range("M137").Value = -4688456.699999999
$ws.Range("N137").Value = -23085573
$ws.Range("H138").Value = 1792.2858
$ws.Range("I138").Value = 1097.3889
$ws.Range("J138").Value = 3043.1
$ws.Range("K138").Value = 3292.1667
$ws.Range("L138").Value = 9129.299999999999
$ws.Range("M138").Value = 1847.8333
$ws.Range("N138").Value = -19409.3
$ws.Range("H141").Value = 1507.3918
$ws.Range("I141").Value = 998.8222
$ws.Range("J141").Value = 2296.5518
$ws.Range("K141").Value = 2996.4666
$ws.Range("L141").Value = 6889.655400000001
$ws.Range("M141").Value = 2183.5334
$ws.Range("N141").Value = -17249.6554

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4747.66
$ws.Range("I32").Value = 3075.889
$ws.Range("J32").Value = 19793.6
$ws.Range("K32").Value = 3075.889
$ws.Range("L32").Value = 19793.6
$ws.Range("M32").Value = -2788.889
$ws.Range("N32").Value = -20367.6
$ws.Range("H132").Value = 108433.234
$ws.Range("I132").Value = 140312.17
$ws.Range("J132").Value = 4102.1816
$ws.Range("K132").Value = 420936.51
$ws.Range("L132").Value = 12306.5448
$ws.Range("M132").Value = -418406.51
$ws.Range("N132").Value = -17366.5448

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 48885.97
$ws.Range("I134").Value = 59968.746
$ws.Range("J134").Value = 1784.1666
$ws.Range("K134").Value = 179906.238
$ws.Range("L134").Value = 5352.4998
$ws.Range("M134").Value = -177371.238
$ws.Range("N134").Value = -10422.4998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1496.4697
$ws.Range("I31").Value = 1215.6364
$ws.Range("J31").Value = 2058.1365
$ws.Range("K31").Value = 1215.6364
$ws.Range("L31").Value = 2058.1365
$ws.Range("M31").Value = -920.6364000000001
$ws.Range("N31").Value = -2648.1365
$ws.Range("H34").Value = 1496.4697
$ws.Range("I34").Value = 1215.6364
$ws.Range("J34").Value = 2058.1365
$ws.Range("K34").Value = 1215.6364
$ws.Range("L34").Value = 2058.1365
$ws.Range("M34").Value = -1013.6364
$ws.Range("N34").Value = -2462.1365
$ws.Range("H58").Value = 1284.449
$ws.Range("I58").Value = 1497.4445
$ws.Range("J58").Value = 694.61536
$ws.Range("K58").Value = 1497.4445
$ws.Range("L58").Value = 694.61536
$ws.Range("M58").Value = -1294.4445
$ws.Range("N58").Value = -1100.61536
$ws.Range("H132").Value = 1915.0962
$ws.Range("I132").Value = 1701.5555
$ws.Range("J132").Value = 3287.8572
$ws.Range("K132").Value = 5104.666499999999
$ws.Range("L132").Value = 9863.571599999999
$ws.Range("M132").Value = -2574.666499999999
$ws.Range("N132").Value = -14923.5716
$ws.Range("H134").Value = 4286.8213
$ws.Range("I134").Value = 4504.4897
$ws.Range("J134").Value = 2763.1428
$ws.Range("K134").Value = 13513.4691
$ws.Range("L134").Value = 8289.428400000001
$ws.Range("M134").Value = -10978.4691
$ws.Range("N134").Value = -13359.4284
$ws.Range("H136").Value = 1284.449
$ws.Range("I136").Value = 1497.4445
$ws.Range("J136").Value = 694.61536
$ws.Range("K136").Value = 4492.333500000001
$ws.Range("L136").Value = 2083.84608
$ws.Range("M136").Value = -1942.333500000001
$ws.Range("N136").Value = -7183.84608

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 87.666664
$ws.Range("I38").Value = 82.55556
$ws.Range("J38").Value = 95.333336
$ws.Range("K38").Value = 247.66668
$ws.Range("L38").Value = 286.000008
$ws.Range("M38").Value = 99.33332000000001
$ws.Range("N38").Value = -980.000008

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 45910
$ws.Range("J133").Value = 45910
$ws.Range("L133").Value = 45910
$ws.Range("N133").Value = -56030

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4764460
$ws.Range("I7").Value = 20001980
$ws.Range("J7").Value = 2734.4375
$ws.Range("K7").Value = 20001980
$ws.Range("L7").Value = 2734.4375
$ws.Range("M7").Value = -20001868
$ws.Range("N7").Value = -2958.4375
$ws.Range("H126").Value = 4764460
$ws.Range("I126").Value = 20001980
$ws.Range("J126").Value = 2734.4375
$ws.Range("K126").Value = 60005940
$ws.Range("L126").Value = 8203.3125
$ws.Range("M126").Value = -60003470
$ws.Range("N126").Value = -13143.3125
$ws.Range("H136").Value = 1659.625
$ws.Range("I136").Value = 1456.2646
$ws.Range("J136").Value = 2153.5
$ws.Range("K136").Value = 4368.793799999999
$ws.Range("L136").Value = 6460.5
$ws.Range("M136").Value = -1818.793799999999
$ws.Range("N136").Value = -11560.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 18237.25
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").ClearContents()
$ws.Range("H42").Value = 17988
$ws.Range("J42").Value = 17988
$ws.Range("L42").Value = 17988
$ws.Range("N42").Value = -18744
$ws.Range("H132").Value = 3958.2856
$ws.Range("I132").Value = 4107.95
$ws.Range("J132").Value = 965
$ws.Range("K132").Value = 12323.85
$ws.Range("L132").Value = 2895
$ws.Range("M132").Value = -9793.849999999999
$ws.Range("N132").Value = -7955
$ws.Range("H136").Value = 1651.3889
$ws.Range("I136").Value = 1668.1111
$ws.Range("J136").Value = 1567.7778
$ws.Range("K136").Value = 5004.3333
$ws.Range("L136").Value = 4703.3334
$ws.Range("M136").Value = -2454.3333
$ws.Range("N136").Value = -9803.3334
